$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.900485992431641
$ws.Range("B1").Value = 4.908373355865479
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 3.920625925064087
$ws.Range("E1").Value = 2.033992290496826
